# VRScreenTrialStructure: add new trial-structure rows (5-11) describing
# distance-calculation changes (2D/3D redo passes, combined-stim experiment).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(5, 11).Value2 = 'first 3D'
$ws.Cells.Item(5, 1).Value2 = '[3]'
$ws.Cells.Item(5, 2).Value2 = '[(0.01, 0.01, 0.01), (0.02, 0.02, 0.02), (0.04, 0.04, 0.04), (0.06, 0.06, 0.06),(0.08, 0.08, 0.08)]'
$ws.Cells.Item(6, 11).Value2 = 'redo 2D with fixes'
$ws.Cells.Item(7, 11).Value2 = 'redo 3D with fixes'
$ws.Cells.Item(6, 5).Value2 = '[0.02, 0.07, 0.15, 0.3]'
$ws.Cells.Item(8, 1).Value2 = '[0, 3]'
$ws.Cells.Item(8, 2).Value2 = '[(0, 0.01, 0.01), (0, 0.02, 0.02), (0, 0.04, 0.04), (0, 0.08, 0.08)]'
$ws.Cells.Item(8, 4).Value2 = '[(0, 0, 0, 1), (0.25, 0.25, 0.25, 1),  (0.5, 0.5, 0.5, 1)]'
$ws.Cells.Item(9, 2).Value2 = '[(0.01, 0.01, 0.01), (0.02, 0.02, 0.02), (0.04, 0.04, 0.04), (0.08, 0.08, 0.08)]'
$ws.Cells.Item(8, 11).Value2 = 'new combined stim, first exp'
$ws.Cells.Item(10, 5).Value2 = '[0.02, 0.07, 0.15]'
$ws.Cells.Item(10, 2).Value2 = '[(0.01, 0.01, 0.01), (0.04, 0.04, 0.04), (0.08, 0.08, 0.08)]'
$ws.Cells.Item(10, 4).Value2 = '[(0, 0, 0, 1), (0.25, 0.25, 0.25, 1), (1, 1, 1, 1)]'
$ws.Cells.Item(11, 5).Value2 = '[0.02, 0.15, 0.30]'
$ws.Cells.Item(11, 2).Value2 = '[(0.02, 0.02, 0.02), (0.06, 0.06, 0.06),(0.08, 0.08, 0.08)]'
$ws.Cells.Item(5, 3).Value2 = '[(1, 1, 1, 1)]'
$ws.Cells.Item(5, 4).Value2 = '[(0, 0, 0, 1), (0.25, 0.25, 0.25, 1)]'
$ws.Cells.Item(5, 5).Value2 = '[0.02, 0.07, 0.15, 0.3, 0.5]'
$ws.Cells.Item(5, 6).Value2 = '[10]'
$ws.Cells.Item(5, 7).Value2 = '[2]'
$ws.Cells.Item(5, 9).Value2 = 5
$ws.Cells.Item(5, 10).Value2 = 3
$ws.Cells.Item(6, 1).Value2 = '[0]'
$ws.Cells.Item(6, 2).Value2 = '[(0.01, 0, 0.01), (0.02, 0, 0.02), (0.04, 0, 0.04), (0.06, 0, 0.06),(0.08, 0, 0.08)]'
$ws.Cells.Item(6, 3).Value2 = '[(1, 1, 1, 1)]'
$ws.Cells.Item(6, 4).Value2 = '[(0, 0, 0, 1), (0.25, 0.25, 0.25, 1)]'
$ws.Cells.Item(6, 6).Value2 = '[10]'
$ws.Cells.Item(6, 7).Value2 = '[0]'
$ws.Cells.Item(6, 9).Value2 = 5
$ws.Cells.Item(6, 10).Value2 = 3
$ws.Cells.Item(7, 1).Value2 = '[3]'
$ws.Cells.Item(7, 2).Value2 = '[(0.01, 0.01, 0.01), (0.02, 0.02, 0.02), (0.04, 0.04, 0.04), (0.06, 0.06, 0.06),(0.08, 0.08, 0.08)]'
$ws.Cells.Item(7, 3).Value2 = '[(1, 1, 1, 1)]'
$ws.Cells.Item(7, 4).Value2 = '[(0, 0, 0, 1), (0.25, 0.25, 0.25, 1)]'
$ws.Cells.Item(7, 5).Value2 = '[0.02, 0.07, 0.15, 0.3]'
$ws.Cells.Item(7, 6).Value2 = '[10]'
$ws.Cells.Item(7, 7).Value2 = '[2]'
$ws.Cells.Item(7, 9).Value2 = 3
$ws.Cells.Item(7, 10).Value2 = 3
$ws.Cells.Item(8, 3).Value2 = '[(1, 1, 1, 1)]'
$ws.Cells.Item(8, 5).Value2 = '[0.02, 0.07, 0.15, 0.3]'
$ws.Cells.Item(8, 6).Value2 = '[10]'
$ws.Cells.Item(8, 7).Value2 = '[0]'
$ws.Cells.Item(8, 9).Value2 = 3
$ws.Cells.Item(8, 10).Value2 = 2
$ws.Cells.Item(9, 1).Value2 = '[0, 3]'
$ws.Cells.Item(9, 3).Value2 = '[(1, 1, 1, 1)]'
$ws.Cells.Item(9, 4).Value2 = '[(0, 0, 0, 1), (0.25, 0.25, 0.25, 1)]'
$ws.Cells.Item(9, 5).Value2 = '[0.02, 0.07, 0.15, 0.3]'
$ws.Cells.Item(9, 6).Value2 = '[10]'
$ws.Cells.Item(9, 7).Value2 = '[0]'
$ws.Cells.Item(9, 9).Value2 = 1
$ws.Cells.Item(9, 10).Value2 = 2
$ws.Cells.Item(10, 1).Value2 = '[0, 3]'
$ws.Cells.Item(10, 3).Value2 = '[(1, 1, 1, 1)]'
$ws.Cells.Item(10, 6).Value2 = '[10]'
$ws.Cells.Item(10, 7).Value2 = '[0]'
$ws.Cells.Item(10, 9).Value2 = 1
$ws.Cells.Item(10, 10).Value2 = 2
$ws.Cells.Item(11, 1).Value2 = '[0, 3]'
$ws.Cells.Item(11, 3).Value2 = '[(1, 1, 1, 1)]'
$ws.Cells.Item(11, 4).Value2 = '[(0, 0, 0, 1), (0.25, 0.25, 0.25, 1), (1, 1, 1, 1)]'
$ws.Cells.Item(11, 6).Value2 = '[10]'
$ws.Cells.Item(11, 7).Value2 = '[0]'
$ws.Cells.Item(11, 9).Value2 = 2
$ws.Cells.Item(11, 10).Value2 = 1

# Widen columns B and D slightly to fit the new, longer entries.
$ws.Columns.Item(2).ColumnWidth = 76.333333333333333
$ws.Columns.Item(4).ColumnWidth = 34.666666666666667

# Leave the selection where the author ended up after entering the data.
$ws.Range("I14").Select()
